$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.051296426740172
$ws.Range("D2").Value = 1.052279358592197
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.058009424959225
$ws.Range("I2").Value = 1.045847733545185
$ws.Range("J2").Value = 1.056325047302772
$ws.Range("K2").Value = 1.055028530387399
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.06074283115631
$ws.Range("N2").Value = 1.057825150032118
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.052553867199837
$ws.Range("D3").Value = 1.053266012509266
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.059456652714152
$ws.Range("I3").Value = 1.046256549744352
$ws.Range("J3").Value = 1.057231013098388
$ws.Range("K3").Value = 1.055827552028815
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.06200240144892
$ws.Range("N3").Value = 1.05873240240308
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.053366433331574
$ws.Range("D4").Value = 1.053903446821404
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.060392445772775
$ws.Range("I4").Value = 1.046519315599123
$ws.Range("J4").Value = 1.057815645464979
$ws.Range("K4").Value = 1.056342963016663
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.062816250216769
$ws.Range("N4").Value = 1.05931786501469
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.053707781027166
$ws.Range("D5").Value = 1.05417118795844
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.060785700037161
$ws.Range("I5").Value = 1.046629361087665
$ws.Range("J5").Value = 1.058061047481508
$ws.Range("K5").Value = 1.056559258821743
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.063158114688741
$ws.Range("N5").Value = 1.059563615530231
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.053765079961031
$ws.Range("D6").Value = 1.054216129088554
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.060851720320146
$ws.Range("I6").Value = 1.046647813546162
$ws.Range("J6").Value = 1.058102229507153
$ws.Range("K6").Value = 1.056595553448836
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.063215499087742
$ws.Range("N6").Value = 1.059604856039075
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.053370995433722
$ws.Range("D7").Value = 1.05390702531693
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.060397701051153
$ws.Range("I7").Value = 1.046520787685919
$ws.Range("J7").Value = 1.057818926018357
$ws.Range("K7").Value = 1.056345854673738
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.062820819311778
$ws.Range("N7").Value = 1.05932115022683
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.051721611009445
$ws.Range("D8").Value = 1.052613010720724
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.058498661972866
$ws.Range("I8").Value = 1.045986261311723
$ws.Range("J8").Value = 1.056631553173869
$ws.Range("K8").Value = 1.055298898144305
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.061168755241634
$ws.Range("N8").Value = 1.058132091176709
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.048806722880665
$ws.Range("D9").Value = 1.050325050553344
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.055147037787856
$ws.Range("I9").Value = 1.045030772217265
$ws.Range("J9").Value = 1.054526963400466
$ws.Range("K9").Value = 1.053441595538103
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.058248394990772
$ws.Range("N9").Value = 1.056024512644428
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.046857526659786
$ws.Range("D10").Value = 1.048794386521896
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.052908766571191
$ws.Range("I10").Value = 1.04438455568911
$ws.Range("J10").Value = 1.053115466363905
$ws.Range("K10").Value = 1.052194886561291
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.056295020329252
$ws.Range("N10").Value = 1.054611011120174
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.046012038975423
$ws.Range("D11").Value = 1.048130286666547
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.051938586746439
$ws.Range("I11").Value = 1.044102528615657
$ws.Range("J11").Value = 1.052502231514275
$ws.Range("K11").Value = 1.051652996159128
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.055447591859822
$ws.Range("N11").Value = 1.053996905406718
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.04569776100969
$ws.Range("D12").Value = 1.047883410167977
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.051578063426907
$ws.Range("I12").Value = 1.043997437189466
$ws.Range("J12").Value = 1.052274137670052
$ws.Range("K12").Value = 1.051451401856748
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.055132572762885
$ws.Range("N12").Value = 1.053768487643074
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.045765185015175
$ws.Range("D13").Value = 1.047936375126348
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.051655404029623
$ws.Range("I13").Value = 1.044019994801437
$ws.Range("J13").Value = 1.052323078709543
$ws.Range("K13").Value = 1.051494658689738
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.055200156667059
$ws.Range("N13").Value = 1.053817498184455
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.045986065303562
$ws.Range("D14").Value = 1.048109883881053
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.051908788980323
$ws.Range("I14").Value = 1.044093848549983
$ws.Range("J14").Value = 1.052483383573148
$ws.Range("K14").Value = 1.051636338695014
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.055421557341451
$ws.Range("N14").Value = 1.053978030699353
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.046122126818751
$ws.Range("D15").Value = 1.048216761724161
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.052064887046039
$ws.Range("I15").Value = 1.044139307966976
$ws.Range("J15").Value = 1.052582111358348
$ws.Range("K15").Value = 1.051723590978356
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.055557936787559
$ws.Range("N15").Value = 1.054076898689332
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.046913607442641
$ws.Range("D16").Value = 1.04883843273128
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.052973132760471
$ws.Range("I16").Value = 1.044403226157027
$ws.Range("J16").Value = 1.053156121305783
$ws.Range("K16").Value = 1.052230806443565
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.056351227104122
$ws.Range("N16").Value = 1.054651723796732
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.0474096845408
$ws.Range("D17").Value = 1.049228037030163
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.053542581016993
$ws.Range("I17").Value = 1.044568181816264
$ws.Range("J17").Value = 1.053515631969891
$ws.Range("K17").Value = 1.052548416417823
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.056848403830649
$ws.Range("N17").Value = 1.05501174500722
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.047698896144131
$ws.Range("D18").Value = 1.049455160226331
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.053874635092694
$ws.Range("I18").Value = 1.044664184451557
$ws.Range("J18").Value = 1.053725131161774
$ws.Range("K18").Value = 1.052733474658626
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.057138244174004
$ws.Range("N18").Value = 1.055221541711986
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.047797485927297
$ws.Range("D19").Value = 1.049532582020159
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.053987840942427
$ws.Range("I19").Value = 1.044696882760097
$ws.Range("J19").Value = 1.053796531613548
$ws.Range("K19").Value = 1.052796541194493
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.057237046209413
$ws.Range("N19").Value = 1.055293043560591
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.047356474864446
$ws.Range("D20").Value = 1.049186249273676
$ws.Range("E20").Value = 0.9894336180360677
$ws.Range("F20").Value = 1.053481494527393
$ws.Range("I20").Value = 1.04455050570229
$ws.Range("J20").Value = 1.053477080311914
$ws.Range("K20").Value = 1.05251436040559
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.056795077487043
$ws.Range("N20").Value = 1.054973138601467
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.045921027847614
$ws.Range("D21").Value = 1.048058795430057
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.051834177764615
$ws.Range("I21").Value = 1.04407210969695
$ws.Range("J21").Value = 1.052436186403124
$ws.Range("K21").Value = 1.051594626103937
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.055356367207319
$ws.Range("N21").Value = 1.053930766503935
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.045017194189809
$ws.Range("D22").Value = 1.047348760354593
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.050797542588367
$ws.Range("I22").Value = 1.043769389726015
$ws.Range("J22").Value = 1.051779933299098
$ws.Range("K22").Value = 1.051014545543987
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.054450364480502
$ws.Range("N22").Value = 1.053273581445227
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.045496459441851
$ws.Range("D23").Value = 1.047725274385045
$ws.Range("E23").Value = 0.9879432794636464
$ws.Range("F23").Value = 1.051347170013926
$ws.Range("I23").Value = 1.043930051188741
$ws.Range("J23").Value = 1.052127997459257
$ws.Range("K23").Value = 1.051322229526344
$ws.Range("L23").Value = 0.9917760702887611
$ws.Range("M23").Value = 1.054930790879557
$ws.Range("N23").Value = 1.053622139896417
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.047380518464606
$ws.Range("D24").Value = 1.049205131752544
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.053509097177189
$ws.Range("I24").Value = 1.044558493436033
$ws.Range("J24").Value = 1.053494500756154
$ws.Range("K24").Value = 1.05252974946412
$ws.Range("L24").Value = 0.9929938892766441
$ws.Range("M24").Value = 1.05681917384315
$ws.Range("N24").Value = 1.054990583784737
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.049561319578997
$ws.Range("D25").Value = 1.050917476093795
$ws.Range("E25").Value = 0.9912096547607051
$ws.Range("F25").Value = 1.056014169738531
$ws.Range("I25").Value = 1.045279408189565
$ws.Range("J25").Value = 1.05507252444058
$ws.Range("K25").Value = 1.0539232412493
$ws.Range("L25").Value = 0.9944092447426416
$ws.Range("M25").Value = 1.059004497939053
$ws.Range("N25").Value = 1.056570848443795
